# "fixed typo in cleaning code"
#
# The States lookup table was missing two territories that the cleaning
# script expects to find ("U.S. Virgin Islands" and "Puerto Rico"), both
# belonging to a "Noncontinental" region. Add them as rows 53-54, plus a
# trailing blank styled row (55) that mirrors the rest of the table's
# formatting. Also a number of rows that were previously force-wrapped to
# two lines (ht=30) no longer need the explicit height now that the table
# was revisited, so let those rows return to the default auto height.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows whose forced 30pt row height is no longer needed -> back to default.
$rowsToAutofit = @(21, 22, 23, 24, 25, 26, 27, 35, 36, 37, 38, 39, 40, 41, 42, 51, 52)
foreach ($r in $rowsToAutofit) {
    $ws.Rows.Item($r).AutoFit()
}

# New territory rows.
$ws.Range("A53").Value = "U.S. Virgin Islands"
$ws.Range("B53").Value = "Noncontinental"
$ws.Range("A54").Value = "Puerto Rico"
$ws.Range("B54").Value = "Noncontinental"

# Match the wrap-text formatting used by the rest of the table and give the
# two new rows the same explicit 30pt height.
$ws.Range("A53:B54").WrapText = $true
$ws.Range("A53:B54").RowHeight = 30

# Trailing blank row keeps the wrap-text cell style but has no value.
$ws.Range("A55").WrapText = $true

# Leave the selection where the user ended up after typing the new data.
$ws.Range("B54").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
